$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range("I2").Value = 3707
$ws.Range("I3").Value = 3819
$ws.Range("H4").Value = 1667
$ws.Range("I4").Value = 902
$ws.Range("I5").Value = 352
$ws.Range("I6").Value = 4307
$ws.Range("H7").Value = 25977
$ws.Range("I7").Value = 13087

$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Range("I2").Value = 134
$ws.Range("I3").Value = 130
$ws.Range("I6").Value = 117
$ws.Range("I7").Value = 422

$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Range("I2").Value = 122
$ws.Range("I3").Value = 177
$ws.Range("I6").Value = 167
$ws.Range("I7").Value = 503

$ws = $wb.Worksheets.Item('Gage Park')
$ws.Range("I6").Value = 39
$ws.Range("I7").Value = 123

$ws = $wb.Worksheets.Item('South Deering')
$ws.Range("I6").Value = 31
$ws.Range("I7").Value = 115

$ws = $wb.Worksheets.Item('New City')
$ws.Range("I6").Value = 89
$ws.Range("I7").Value = 286

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range("I5").Value = 40
$ws.Range("I8").Value = 795
$ws.Range("I11").Value = 205
$ws.Range("I15").Value = 157
$ws.Range("I19").Value = 350
$ws.Range("I20").Value = 321
$ws.Range("I21").Value = 73
$ws.Range("I27").Value = 122
$ws.Range("I29").Value = 843
$ws.Range("I31").Value = 123
$ws.Range("I33").Value = 589
$ws.Range("I34").Value = 63
$ws.Range("I37").Value = 422
$ws.Range("I42").Value = 458
$ws.Range("I43").Value = 115
$ws.Range("I45").Value = 25
$ws.Range("I48").Value = 176
$ws.Range("I49").Value = 107
$ws.Range("I51").Value = 126
$ws.Range("I52").Value = 288
$ws.Range("I53").Value = 143
$ws.Range("I54").Value = 298
$ws.Range("H63").Value = 207
$ws.Range("I65").Value = 286
$ws.Range("I66").Value = 34
$ws.Range("I67").Value = 503
$ws.Range("I76").Value = 199
$ws.Range("I78").Value = 187
$ws.Range("I79").Value = 352
$ws.Range("I83").Value = 265
$ws.Range("I84").Value = 115
$ws.Range("I85").Value = 592
$ws.Range("I86").Value = 79
$ws.Range("I87").Value = 22
$ws.Range("I90").Value = 165
$ws.Range("I94").Value = 121
$ws.Range("I95").Value = 208
$ws.Range("I97").Value = 95
$ws.Range("I98").Value = 87
$ws.Range("H101").Value = 25977
$ws.Range("I101").Value = 13087

$ws = $wb.Worksheets.Item('South Chicago')
$ws.Range("I2").Value = 93
$ws.Range("I7").Value = 265

$ws = $wb.Worksheets.Item('West Pullman')
$ws.Range("I2").Value = 74
$ws.Range("I3").Value = 78
$ws.Range("I7").Value = 208

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range("I2").Value = 140
$ws.Range("I3").Value = 214
$ws.Range("I5").Value = 20
$ws.Range("I7").Value = 589

$ws = $wb.Worksheets.Item('Lincoln Park')
$ws.Range("I2").Value = 20
$ws.Range("I6").Value = 68
$ws.Range("I7").Value = 107

$ws = $wb.Worksheets.Item('Loop')
$ws.Range("I2").Value = 68
$ws.Range("I4").Value = 21
$ws.Range("I6").Value = 150
$ws.Range("I7").Value = 298

$ws = $wb.Worksheets.Item('Englewood')
$ws.Range("I2").Value = 253
$ws.Range("I3").Value = 287
$ws.Range("I7").Value = 843

$ws = $wb.Worksheets.Item('Chatham')
$ws.Range("I2").Value = 133
$ws.Range("I4").Value = 15
$ws.Range("I7").Value = 350

$ws = $wb.Worksheets.Item('Lake View')
$ws.Range("I6").Value = 103
$ws.Range("I7").Value = 176

$ws = $wb.Worksheets.Item('River North')
$ws.Range("I6").Value = 84
$ws.Range("I7").Value = 199

$ws = $wb.Worksheets.Item('South Shore')
$ws.Range("I3").Value = 238
$ws.Range("I7").Value = 592

$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Range("I2").Value = 125
$ws.Range("I3").Value = 156
$ws.Range("I7").Value = 458

$ws = $wb.Worksheets.Item('Rogers Park')
$ws.Range("I2").Value = 38
$ws.Range("I6").Value = 76
$ws.Range("I7").Value = 187

$ws = $wb.Worksheets.Item('Chinatown')
$ws.Range("I6").Value = 56
$ws.Range("I7").Value = 73

$ws = $wb.Worksheets.Item('Roseland')
$ws.Range("I2").Value = 104
$ws.Range("I4").Value = 21
$ws.Range("I7").Value = 352

$ws = $wb.Worksheets.Item('Chicago Lawn')
$ws.Range("I2").Value = 90
$ws.Range("I4").Value = 22
$ws.Range("I6").Value = 101
$ws.Range("I7").Value = 321

$ws = $wb.Worksheets.Item('Little Village')
$ws.Range("I2").Value = 80
$ws.Range("I3").Value = 96
$ws.Range("I4").Value = 30
$ws.Range("I6").Value = 72
$ws.Range("I7").Value = 288

$ws = $wb.Worksheets.Item('Garfield Ridge')
$ws.Range("I3").Value = 23
$ws.Range("I7").Value = 63

$ws = $wb.Worksheets.Item('West Loop')
$ws.Range("I2").Value = 22
$ws.Range("I7").Value = 121

$ws = $wb.Worksheets.Item('Brighton Park')
$ws.Range("I2").Value = 50
$ws.Range("I4").Value = 12
$ws.Range("I7").Value = 157

$ws = $wb.Worksheets.Item('Wicker Park')
$ws.Range("I6").Value = 53
$ws.Range("I7").Value = 87

$ws = $wb.Worksheets.Item('North Center')
$ws.Range("I6").Value = 15
$ws.Range("I7").Value = 34

$ws = $wb.Worksheets.Item('Belmont Cragin')
$ws.Range("I6").Value = 48
$ws.Range("I7").Value = 205

$ws = $wb.Worksheets.Item('West Town')
$ws.Range("I6").Value = 55
$ws.Range("I7").Value = 95

$ws = $wb.Worksheets.Item('Austin')
$ws.Range("I2").Value = 248
$ws.Range("I3").Value = 218
$ws.Range("I6").Value = 259
$ws.Range("I7").Value = 795

$ws = $wb.Worksheets.Item('Armour Square')
$ws.Range("I3").Value = 12
$ws.Range("I7").Value = 40

$ws = $wb.Worksheets.Item('Edgewater')
$ws.Range("I6").Value = 50
$ws.Range("I7").Value = 122

$ws = $wb.Worksheets.Item('Streeterville')
$ws.Range("I3").Value = 6
$ws.Range("I7").Value = 79

$ws = $wb.Worksheets.Item('Washington Heights')
$ws.Range("I2").Value = 54
$ws.Range("I7").Value = 165

$ws = $wb.Worksheets.Item('Little Italy, UIC')
$ws.Range("I4").Value = 12
$ws.Range("I7").Value = 126

$ws = $wb.Worksheets.Item('Hyde Park')
$ws.Range("I6").Value = 67
$ws.Range("I7").Value = 115

$ws = $wb.Worksheets.Item('Logan Square')
$ws.Range("I4").Value = 10
$ws.Range("I7").Value = 143

$ws = $wb.Worksheets.Item('Jackson Park')
$ws.Range("I6").Value = 8
$ws.Range("I7").Value = 25

$ws = $wb.Worksheets.Item('Ukrainian Village')
$ws.Range("I6").Value = 11
$ws.Range("I7").Value = 22
